$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.206.63"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "'1.906.16"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'306.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.5248"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.80%  "

$ws.Range("D8").Value = "'0.3772"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.32%  "

$ws.Range("D9").Value = "'0.07252"
$ws.Range("D9").ClearFormats()

$ws.Range("D10").Value = "'21.16"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").Value = "'0.9010"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "'0.08419"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +10.39%  "

$ws.Range("D13").Value = "'1.904.88"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("D14").Value = "'94.95"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").Value = "'5.287"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").Value = "'0.000008622"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.27%  "

$ws.Range("D18").Value = "'14.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "'27.238.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("D21").Value = "'5.068"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").Value = "'2.145.77"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("D24").Value = "'6.437"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").Value = "'147.53"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").Value = "'2.289"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.23%  "

$ws.Range("E27").Value = "  -2.08%  "

$ws.Range("D28").Value = "'18.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("D29").Value = "'114.96"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.821"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.920"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("D32").Value = "'0.09283"
$ws.Range("D32").ClearFormats()

$ws.Range("D33").Value = "'0.8105"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.18%  "

$ws.Range("D34").Value = "'0.05068"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "'1.243"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.88%  "

$ws.Range("D36").Value = "'2.955"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").Value = "'3.381"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.29%  "

$ws.Range("D38").Value = "'2.610"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").Value = "'0.5750"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.41%  "

$ws.Range("D40").Value = "'0.01991"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").Value = "'1.076"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "'6.661"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.29%  "

$ws.Range("D43").Value = "'8.994"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'117.51"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "'0.1514"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D46").Value = "'0.4857"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.16%  "

$ws.Range("D47").Value = "'10.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").Value = "'1.619"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.23%  "

$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("D51").Value = "'63.98"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.35%  "

